# Update the cryptos price/volume snapshot (columns D "Price" and E
# "Volume(1h)") for the rows whose figures changed in the latest refresh.
#
# Columns D/E are stored as plain text in the workbook (e.g. "64.012.23",
# "  -0.54%  "), not numbers, so prices that also happen to parse as valid
# numeric literals (e.g. "571.78", "0.120", "2.00") must be written with a
# Text number format first - otherwise Excel's COM layer would silently
# coerce them to numbers and mangle formatting-significant digits such as
# trailing zeros (e.g. "0.120" -> 0.12) or thousands-style dot grouping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text would otherwise be auto-converted to a number.
$textCells = @(
    "D5","D6","D10","D11","D15","D16","D19","D20","D22","D24","D26","D27",
    "D29","D30","D31","D32","D35","D37","D38","D39","D41","D43","D44","D45",
    "D47","D49"
)
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "64.015.69"
$ws.Range("E2").Value = "  -0.42%  "

$ws.Range("D3").Value = "3.432.58"
$ws.Range("E3").Value = "  +1.03%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "571.78"
$ws.Range("E5").Value = "  +0.61%  "

$ws.Range("D6").Value = "159.68"
$ws.Range("E6").Value = "  +2.22%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "3.433.53"
$ws.Range("E8").Value = "  +0.77%  "

$ws.Range("E9").Value = "  -9.63%  "

$ws.Range("D10").Value = "7.27"
$ws.Range("E10").Value = "  +1.96%  "

$ws.Range("D11").Value = "0.120"
$ws.Range("E11").Value = "  -1.79%  "

$ws.Range("E12").Value = "  -3.47%  "

$ws.Range("D13").Value = "4.021.65"
$ws.Range("E13").Value = "  +0.93%  "

$ws.Range("E14").Value = "  +0.51%  "

$ws.Range("D15").Value = "27.19"
$ws.Range("E15").Value = "  -0.87%  "

$ws.Range("D16").Value = "0.0000174"
$ws.Range("E16").Value = "  -6.82%  "

$ws.Range("D17").Value = "64.055.03"
$ws.Range("E17").Value = "  -0.35%  "

$ws.Range("D18").Value = "3.437.62"
$ws.Range("E18").Value = "  -0.35%  "

$ws.Range("D19").Value = "6.11"
$ws.Range("E19").Value = "  -2.70%  "

$ws.Range("D20").Value = "13.63"
$ws.Range("E20").Value = "  -0.84%  "

$ws.Range("E21").Value = "  +1.98%  "

$ws.Range("D22").Value = "7.86"
$ws.Range("E22").Value = "  -1.43%  "

$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("D24").Value = "71.40"
$ws.Range("E24").Value = "  -0.46%  "

$ws.Range("E25").Value = "  -5.00%  "

$ws.Range("D26").Value = "0.0000117"
$ws.Range("E26").Value = "  -0.87%  "

$ws.Range("D27").Value = "9.69"
$ws.Range("E27").Value = "  -5.52%  "

$ws.Range("E28").Value = "  +0.31%  "

$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.21%  "

$ws.Range("D30").Value = "6.06"
$ws.Range("E30").Value = "  -1.12%  "

$ws.Range("D31").Value = "1.40"
$ws.Range("E31").Value = "  -4.42%  "

$ws.Range("D32").Value = "2.00"
$ws.Range("E32").Value = "  -0.04%  "

$ws.Range("E33").Value = "  +0.14%  "

$ws.Range("D35").Value = "6.96"
$ws.Range("E35").Value = "  -2.19%  "

$ws.Range("E36").Value = "  -4.53%  "

$ws.Range("D37").Value = "160.86"
$ws.Range("E37").Value = "  +0.87%  "

$ws.Range("D38").Value = "0.849"
$ws.Range("E38").Value = "  +10.84%  "

$ws.Range("D39").Value = "1.85"
$ws.Range("E39").Value = "  -2.06%  "

$ws.Range("D40").Value = "2.836.05"
$ws.Range("E40").Value = "  -1.35%  "

$ws.Range("D41").Value = "26.12"
$ws.Range("E41").Value = "  -0.19%  "

$ws.Range("E42").Value = "  -4.24%  "

$ws.Range("D43").Value = "43.06"
$ws.Range("E43").Value = "  +0.39%  "

$ws.Range("D44").Value = "26.45"
$ws.Range("E44").Value = "  +2.20%  "

$ws.Range("D45").Value = "6.43"
$ws.Range("E45").Value = "  -7.19%  "

$ws.Range("E46").Value = "  -5.06%  "

$ws.Range("D47").Value = "0.0305"
$ws.Range("E47").Value = "  -2.79%  "

$ws.Range("E48").Value = "  +11.04%  "

$ws.Range("D49").Value = "333.87"
$ws.Range("E49").Value = "  +3.98%  "

$ws.Range("E50").Value = "  -1.97%  "

$ws.Range("E51").Value = "  -6.08%  "
